$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8078123927116394
$ws.Range("B1").Value = 0.9029212594032288
$ws.Range("C1").Value = 1.096034526824951
$ws.Range("D1").Value = 1.306583642959595
$ws.Range("E1").Value = 1.621540904045105
